$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top, shifting existing data down
$ws.Rows.Item(1).Insert()

# Add the heading text in the newly inserted row
$ws.Range("A1").Value = "Heading"
$ws.Range("A1").Font.Bold = $true
